# Auto-generated edit script for cs-en-us-106pct.xlsx weekly update
# (New crime data collected: shifts reporting week from 12/25/2023-12/31/2023
#  to 1/1/2024-1/7/2024, Volume 30 No. 52 -> Volume 31 No. 1, and refreshes all
#  weekly/28-day/YTD/2-year/historical crime-complaint figures accordingly.)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Masthead / header text changes (Volume, report week, year-span labels) ---
$ws.Range('A8').Value = "Volume 31   Number  1"
$ws.Range('C9').Value = "Report Covering the Week  1/1/2024  Through  1/7/2024"
$ws.Range('M12').Value = "14 Year (2010)"
$ws.Range('N12').Value = "31 Year (1993)"
$ws.Range('K35').Value = "''23 vs '01"
$ws.Range('L35').Value = "''23 vs '98"
$ws.Range('M35').Value = "''23 vs '93"
$ws.Range('N35').Value = "''23 vs '90"

# --- Year column headers on the crime-complaints tables (2023/2022 -> 2024/2023) ---
# --- and refreshed numeric figures for Week to Date / 28 Day / Year to Date / 2 Year ---
# --- and Historical Perspective sections (style/type unchanged from prior week) ---
$ws.Range('C13').Value = 2024
$ws.Range('D13').Value = 2023
$ws.Range('F13').Value = 2024
$ws.Range('G13').Value = 2023
$ws.Range('I13').Value = 2024
$ws.Range('J13').Value = 2023
$ws.Range('C16').Value = 1
$ws.Range('D16').Value = 3
$ws.Range('E16').Value = -66.666666666666
$ws.Range('F16').Value = 13
$ws.Range('G16').Value = 8
$ws.Range('H16').Value = 62.5
$ws.Range('I16').Value = 1
$ws.Range('J16').Value = 3
$ws.Range('K16').Value = -66.666666666666
$ws.Range('L16').Value = 0
$ws.Range('M16').Value = -92.307692307692
$ws.Range('N16').Value = -95.652173913043
$ws.Range('C17').Value = 6
$ws.Range('D17').Value = 10
$ws.Range('E17').Value = -40
$ws.Range('F17').Value = 33
$ws.Range('G17').Value = 29
$ws.Range('H17').Value = 13.793103448275
$ws.Range('I17').Value = 6
$ws.Range('J17').Value = 10
$ws.Range('K17').Value = -40
$ws.Range('L17').Value = 50
$ws.Range('M17').Value = 50
$ws.Range('N17').Value = -25
$ws.Range('C18').Value = 9
$ws.Range('D18').Value = 4
$ws.Range('E18').Value = 125
$ws.Range('F18').Value = 23
$ws.Range('G18').Value = 8
$ws.Range('H18').Value = 187.5
$ws.Range('I18').Value = 9
$ws.Range('J18').Value = 4
$ws.Range('K18').Value = 125
$ws.Range('L18').Value = 200
$ws.Range('M18').Value = 80
$ws.Range('N18').Value = -73.529411764705
$ws.Range('C19').Value = 4
$ws.Range('D19').Value = 11
$ws.Range('E19').Value = -63.636363636363
$ws.Range('F19').Value = 39
$ws.Range('G19').Value = 52
$ws.Range('H19').Value = -25
$ws.Range('I19').Value = 4
$ws.Range('J19').Value = 11
$ws.Range('K19').Value = -63.636363636363
$ws.Range('L19').Value = -66.666666666666
$ws.Range('M19').Value = 33.333333333333
$ws.Range('N19').Value = -69.230769230769
$ws.Range('C20').Value = 13
$ws.Range('D20').Value = 6
$ws.Range('E20').Value = 116.666666666667
$ws.Range('F20').Value = 33
$ws.Range('G20').Value = 21
$ws.Range('H20').Value = 57.142857142857
$ws.Range('I20').Value = 13
$ws.Range('J20').Value = 6
$ws.Range('K20').Value = 116.666666666667
$ws.Range('L20').Value = 333.333333333333
$ws.Range('M20').Value = 85.714285714285
$ws.Range('N20').Value = -80.882352941176
$ws.Range('C21').Value = 33
$ws.Range('D21').Value = 34
$ws.Range('E21').Value = -2.941176470588
$ws.Range('F21').Value = 142
$ws.Range('G21').Value = 118
$ws.Range('H21').Value = 20.338983050847
$ws.Range('I21').Value = 33
$ws.Range('J21').Value = 34
$ws.Range('K21').Value = -2.941176470588
$ws.Range('L21').Value = 43.478260869565
$ws.Range('M21').Value = 3.125
$ws.Range('N21').Value = -77.397260273972
$ws.Range('L22').Value = -100
$ws.Range('C24').Value = 20
$ws.Range('D24').Value = 13
$ws.Range('E24').Value = 53.846153846153
$ws.Range('F24').Value = 80
$ws.Range('G24').Value = 93
$ws.Range('H24').Value = -13.978494623655
$ws.Range('I24').Value = 20
$ws.Range('J24').Value = 13
$ws.Range('K24').Value = 53.846153846153
$ws.Range('L24').Value = 100
$ws.Range('M24').Value = 122.222222222222
$ws.Range('C25').Value = 11
$ws.Range('D25').Value = 10
$ws.Range('E25').Value = 10
$ws.Range('F25').Value = 34
$ws.Range('G25').Value = 34
$ws.Range('H25').Value = 0
$ws.Range('I25').Value = 11
$ws.Range('J25').Value = 10
$ws.Range('K25').Value = 10
$ws.Range('L25').Value = 120
$ws.Range('M25').Value = 10
$ws.Range('F27').Value = 2
$ws.Range('G27').Value = 3
$ws.Range('H27').Value = -33.333333333333
$ws.Range('J27').Value = 2
$ws.Range('K27').Value = -100
$ws.Range('J30').Value = 1
$ws.Range('K30').Value = -100
$ws.Range('J35').Value = 2023
$ws.Range('J36').Value = 2
$ws.Range('K36').Value = -50
$ws.Range('L36').Value = -71.428571428571
$ws.Range('M36').Value = -90.47619047619
$ws.Range('N36').Value = -86.666666666666
$ws.Range('J37').Value = 22
$ws.Range('K37').Value = 15.78947368421
$ws.Range('L37').Value = -45
$ws.Range('M37').Value = -33.333333333333
$ws.Range('N37').Value = -29.032258064516
$ws.Range('J38').Value = 208
$ws.Range('K38').Value = -39.53488372093
$ws.Range('L38').Value = -61.834862385321
$ws.Range('M38').Value = -76.65544332211
$ws.Range('N38').Value = -75
$ws.Range('J39').Value = 378
$ws.Range('K39').Value = 56.846473029045
$ws.Range('L39').Value = 27.272727272727
$ws.Range('M39').Value = -7.35294117647
$ws.Range('N39').Value = 14.893617021276
$ws.Range('J40').Value = 146
$ws.Range('K40').Value = -64.988009592326
$ws.Range('L40').Value = -77.945619335347
$ws.Range('M40').Value = -88.130081300813
$ws.Range('N40').Value = -89.169139465875
$ws.Range('J41').Value = 616
$ws.Range('K41').Value = 61.679790026246
$ws.Range('L41').Value = 40
$ws.Range('M41').Value = 3.529411764705
$ws.Range('N41').Value = -13.846153846153
$ws.Range('J42').Value = 263
$ws.Range('K42').Value = -66.195372750642
$ws.Range('L42').Value = -80.802919708029
$ws.Range('M42').Value = -92.039951573849
$ws.Range('N42').Value = -91.43601432758
$ws.Range('J43').Value = 1635
$ws.Range('K43').Value = -25.137362637362
$ws.Range('L43').Value = -51.353763760785
$ws.Range('M43').Value = -74.776303609996
$ws.Range('N43').Value = -74.215423434789

# --- Cells whose displayed kind flips between a numeric figure and the '0'/'***.*'
#     placeholder text used when a category has no reportable data this period ---
$ws.Range('C23').Copy()
$ws.Range('C14').PasteSpecial(-4122)
$ws.Range('C23').Copy()
$ws.Range('C14').PasteSpecial(-4163)
$ws.Range('C23').Copy()
$ws.Range('G14').PasteSpecial(-4122)
$ws.Range('C23').Copy()
$ws.Range('G14').PasteSpecial(-4163)
$ws.Range('E23').Copy()
$ws.Range('H14').PasteSpecial(-4122)
$ws.Range('E23').Copy()
$ws.Range('H14').PasteSpecial(-4163)
$ws.Range('C23').Copy()
$ws.Range('I14').PasteSpecial(-4122)
$ws.Range('C23').Copy()
$ws.Range('I14').PasteSpecial(-4163)
$ws.Range('C23').Copy()
$ws.Range('J14').PasteSpecial(-4122)
$ws.Range('C23').Copy()
$ws.Range('J14').PasteSpecial(-4163)
$ws.Range('E23').Copy()
$ws.Range('K14').PasteSpecial(-4122)
$ws.Range('E23').Copy()
$ws.Range('K14').PasteSpecial(-4163)
$ws.Range('E23').Copy()
$ws.Range('L14').PasteSpecial(-4122)
$ws.Range('E23').Copy()
$ws.Range('L14').PasteSpecial(-4163)
$ws.Range('E23').Copy()
$ws.Range('M14').PasteSpecial(-4122)
$ws.Range('E23').Copy()
$ws.Range('M14').PasteSpecial(-4163)
$ws.Range('E23').Copy()
$ws.Range('N14').PasteSpecial(-4122)
$ws.Range('E23').Copy()
$ws.Range('N14').PasteSpecial(-4163)
$ws.Range('C23').Copy()
$ws.Range('I15').PasteSpecial(-4122)
$ws.Range('C23').Copy()
$ws.Range('I15').PasteSpecial(-4163)
$ws.Range('C23').Copy()
$ws.Range('J15').PasteSpecial(-4122)
$ws.Range('C23').Copy()
$ws.Range('J15').PasteSpecial(-4163)
$ws.Range('E23').Copy()
$ws.Range('K15').PasteSpecial(-4122)
$ws.Range('E23').Copy()
$ws.Range('K15').PasteSpecial(-4163)
$ws.Range('E23').Copy()
$ws.Range('L15').PasteSpecial(-4122)
$ws.Range('E23').Copy()
$ws.Range('L15').PasteSpecial(-4163)
$ws.Range('E23').Copy()
$ws.Range('M15').PasteSpecial(-4122)
$ws.Range('E23').Copy()
$ws.Range('M15').PasteSpecial(-4163)
$ws.Range('E23').Copy()
$ws.Range('N15').PasteSpecial(-4122)
$ws.Range('E23').Copy()
$ws.Range('N15').PasteSpecial(-4163)
$ws.Range('C23').Copy()
$ws.Range('D22').PasteSpecial(-4122)
$ws.Range('C23').Copy()
$ws.Range('D22').PasteSpecial(-4163)
$ws.Range('E23').Copy()
$ws.Range('E22').PasteSpecial(-4122)
$ws.Range('E23').Copy()
$ws.Range('E22').PasteSpecial(-4163)
$ws.Range('C23').Copy()
$ws.Range('I22').PasteSpecial(-4122)
$ws.Range('C23').Copy()
$ws.Range('I22').PasteSpecial(-4163)
$ws.Range('C23').Copy()
$ws.Range('J22').PasteSpecial(-4122)
$ws.Range('C23').Copy()
$ws.Range('J22').PasteSpecial(-4163)
$ws.Range('E23').Copy()
$ws.Range('K22').PasteSpecial(-4122)
$ws.Range('E23').Copy()
$ws.Range('K22').PasteSpecial(-4163)
$ws.Range('E23').Copy()
$ws.Range('M22').PasteSpecial(-4122)
$ws.Range('E23').Copy()
$ws.Range('M22').PasteSpecial(-4163)
$ws.Range('C23').Copy()
$ws.Range('I26').PasteSpecial(-4122)
$ws.Range('C23').Copy()
$ws.Range('I26').PasteSpecial(-4163)
$ws.Range('C23').Copy()
$ws.Range('J26').PasteSpecial(-4122)
$ws.Range('C23').Copy()
$ws.Range('J26').PasteSpecial(-4163)
$ws.Range('E23').Copy()
$ws.Range('K26').PasteSpecial(-4122)
$ws.Range('E23').Copy()
$ws.Range('K26').PasteSpecial(-4163)
$ws.Range('E23').Copy()
$ws.Range('L26').PasteSpecial(-4122)
$ws.Range('E23').Copy()
$ws.Range('L26').PasteSpecial(-4163)
$ws.Range('C23').Copy()
$ws.Range('C27').PasteSpecial(-4122)
$ws.Range('C23').Copy()
$ws.Range('C27').PasteSpecial(-4163)
$ws.Range('F14').Copy()
$ws.Range('D27').PasteSpecial(-4122)
$ws.Range('D27').Value = 2
$ws.Range('H22').Copy()
$ws.Range('E27').PasteSpecial(-4122)
$ws.Range('E27').Value = -100
$ws.Range('C23').Copy()
$ws.Range('I27').PasteSpecial(-4122)
$ws.Range('C23').Copy()
$ws.Range('I27').PasteSpecial(-4163)
$ws.Range('E23').Copy()
$ws.Range('L27').PasteSpecial(-4122)
$ws.Range('E23').Copy()
$ws.Range('L27').PasteSpecial(-4163)
$ws.Range('C23').Copy()
$ws.Range('C28').PasteSpecial(-4122)
$ws.Range('C23').Copy()
$ws.Range('C28').PasteSpecial(-4163)
$ws.Range('C23').Copy()
$ws.Range('I28').PasteSpecial(-4122)
$ws.Range('C23').Copy()
$ws.Range('I28').PasteSpecial(-4163)
$ws.Range('C23').Copy()
$ws.Range('J28').PasteSpecial(-4122)
$ws.Range('C23').Copy()
$ws.Range('J28').PasteSpecial(-4163)
$ws.Range('E23').Copy()
$ws.Range('K28').PasteSpecial(-4122)
$ws.Range('E23').Copy()
$ws.Range('K28').PasteSpecial(-4163)
$ws.Range('E23').Copy()
$ws.Range('L28').PasteSpecial(-4122)
$ws.Range('E23').Copy()
$ws.Range('L28').PasteSpecial(-4163)
$ws.Range('E23').Copy()
$ws.Range('M28').PasteSpecial(-4122)
$ws.Range('E23').Copy()
$ws.Range('M28').PasteSpecial(-4163)
$ws.Range('E23').Copy()
$ws.Range('N28').PasteSpecial(-4122)
$ws.Range('E23').Copy()
$ws.Range('N28').PasteSpecial(-4163)
$ws.Range('C23').Copy()
$ws.Range('C29').PasteSpecial(-4122)
$ws.Range('C23').Copy()
$ws.Range('C29').PasteSpecial(-4163)
$ws.Range('C23').Copy()
$ws.Range('I29').PasteSpecial(-4122)
$ws.Range('C23').Copy()
$ws.Range('I29').PasteSpecial(-4163)
$ws.Range('C23').Copy()
$ws.Range('J29').PasteSpecial(-4122)
$ws.Range('C23').Copy()
$ws.Range('J29').PasteSpecial(-4163)
$ws.Range('E23').Copy()
$ws.Range('K29').PasteSpecial(-4122)
$ws.Range('E23').Copy()
$ws.Range('K29').PasteSpecial(-4163)
$ws.Range('E23').Copy()
$ws.Range('L29').PasteSpecial(-4122)
$ws.Range('E23').Copy()
$ws.Range('L29').PasteSpecial(-4163)
$ws.Range('E23').Copy()
$ws.Range('M29').PasteSpecial(-4122)
$ws.Range('E23').Copy()
$ws.Range('M29').PasteSpecial(-4163)
$ws.Range('E23').Copy()
$ws.Range('N29').PasteSpecial(-4122)
$ws.Range('E23').Copy()
$ws.Range('N29').PasteSpecial(-4163)
$ws.Range('F14').Copy()
$ws.Range('D30').PasteSpecial(-4122)
$ws.Range('D30').Value = 1
$ws.Range('H22').Copy()
$ws.Range('E30').PasteSpecial(-4122)
$ws.Range('E30').Value = -100
$ws.Range('F14').Copy()
$ws.Range('G30').PasteSpecial(-4122)
$ws.Range('G30').Value = 1
$ws.Range('H22').Copy()
$ws.Range('H30').PasteSpecial(-4122)
$ws.Range('H30').Value = 0
$ws.Range('C23').Copy()
$ws.Range('I30').PasteSpecial(-4122)
$ws.Range('C23').Copy()
$ws.Range('I30').PasteSpecial(-4163)
$ws.Range('E23').Copy()
$ws.Range('L30').PasteSpecial(-4122)
$ws.Range('E23').Copy()
$ws.Range('L30').PasteSpecial(-4163)

$excel.CutCopyMode = 0
